$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts existing rows 2-17 down to 3-18)
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the 2004 data
$ws.Cells.Item(2, 1).Value = 2004
$ws.Cells.Item(2, 2).Value = "suome energia hiilidioksid polto ilmastonmuutoks"
$ws.Cells.Item(2, 3).Value = "selonteo uhk suome sotilaallis main"
